$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "认证方式" (auth_type) column. This is column P (16th column) on
# the server sheet. Deleting it shifts 用户密码/SSH端口/宿主机 one column left
# (Q->P, R->Q, S->R) and drops the associated data validation + shared
# string automatically.
$ws.Columns.Item(16).Delete()

# Add the 12 new server rows (JR002..JR013), matching the pattern used for
# row 2 (JR001) but with status/date fields like rows 5+ in the template.
$newRows = @(
  @{Row=5;  Asset="JR002"; IP="10.1.19.11"},
  @{Row=6;  Asset="JR003"; IP="10.1.19.12"},
  @{Row=7;  Asset="JR004"; IP="10.1.19.13"},
  @{Row=8;  Asset="JR005"; IP="10.1.19.14"},
  @{Row=9;  Asset="JR006"; IP="10.1.19.15"},
  @{Row=10; Asset="JR007"; IP="10.1.19.16"},
  @{Row=11; Asset="JR008"; IP="10.1.19.17"},
  @{Row=12; Asset="JR009"; IP="10.1.19.18"},
  @{Row=13; Asset="JR010"; IP="10.1.19.19"},
  @{Row=14; Asset="JR011"; IP="10.1.19.20"},
  @{Row=15; Asset="JR012"; IP="10.1.19.21"},
  @{Row=16; Asset="JR013"; IP="10.1.19.22"}
)

# Fill in the non-text-unique columns for every new row first (these reuse
# existing shared strings / are numeric, so ordering doesn't affect the
# shared-string table), then lay down the new, never-seen-before strings
# (asset tag in B, then management IP in F) in the same row-major / column
# order they were originally typed: row 5's B & F together, then every
# remaining row's B down the column, then every remaining row's F down the
# column.
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = "server"
    $ws.Range("D$row").Value = 1
    $ws.Range("E$row").Value = 3
    $ws.Range("G$row").Value = 2
    $ws.Range("H$row").Value = 1
    $ws.Range("I$row").Value = 1

    $ws.Range("J$row").Value = 39668
    $ws.Range("J$row").NumberFormat = "m/d/yy"
    $ws.Range("K$row").Value = 43320
    $ws.Range("K$row").NumberFormat = "m/d/yy"

    $ws.Range("L$row").Value = 800
    $ws.Range("N$row").Value = 0
    $ws.Range("O$row").Value = "root"

    $ws.Range("P$row").Value = "123456"
    $ws.Range("P$row").NumberFormat = "@"

    $ws.Range("Q$row").Value = 22
}

$ws.Range("B5").Value = $newRows[0].Asset
$ws.Range("F5").Value = $newRows[0].IP

for ($i = 1; $i -lt $newRows.Count; $i++) {
    $ws.Range("B$($newRows[$i].Row)").Value = $newRows[$i].Asset
}

for ($i = 1; $i -lt $newRows.Count; $i++) {
    $ws.Range("F$($newRows[$i].Row)").Value = $newRows[$i].IP
}

# Move the active selection to C15 and drop the frozen topLeftCell scroll
# position (now scrolled back to the sheet origin).
$ws.Range("C15").Select()
